# Error Calculations and Plots
# Applies the missing-data edit to the worksheet:
#  - removes the "RM 232" row (row 26) and the "SC 92" row (row 28, original numbering)
#  - fills in / blanks out a number of individual cells to match the new
#    "missing data" pattern used by the imputation exercise

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two data rows (delete bottom-most first so row numbers above
#     stay valid for the second delete) ---
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"

# --- Cell-level value changes on the remaining 32 data rows ---

# Row 5 (RM 14): D column stays, E5 becomes missing
$ws.Range("E5").ClearContents()

# Row 6 (RM 21): F6 gets its value back
$ws.Range("F6").Value = 16.43

# Row 8 (RM 38): E8 gets its value back
$ws.Range("E8").Value = -6.6

# Row 11 (RM 58): F11 gets its value back
$ws.Range("F11").Value = 17.65

# Row 12 (RM 81): E12 and F12 become missing
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

# Row 14 (RM 90): E14 gets its value back
$ws.Range("E14").Value = -5.4

# Row 17 (RM 116): F17 becomes missing
$ws.Range("F17").ClearContents()

# Row 18 (RM 120): E18 becomes missing
$ws.Range("E18").ClearContents()

# Row 25 (RM 145): F25 gets its value back
$ws.Range("F25").Value = 16.6

# Row 26 (now SC 5): C26 gets its value back
$ws.Range("C26").Value = 10.8

# Row 27 (now SC 101): C27 becomes missing
$ws.Range("C27").ClearContents()

# Row 31 (now SC 132): F31 becomes missing
$ws.Range("F31").ClearContents()

# Row 32 (now SC 193): F32 becomes missing
$ws.Range("F32").ClearContents()

# Row 33 (now SC 232): D33 gets its value back
$ws.Range("D33").Value = -14.1
